# The Slide Master and every Slide Layout in this deck carry a "Date
# Placeholder" shape whose text is the cached output of an automatic
# <a:fld type="datetimeFigureOut"> field ("10/8/2023"). Re-point that
# cached text at the following day, "10/9/2023", everywhere it occurs
# (Slide Master + all Custom Layouts).

$p = $ppt.ActivePresentation

$oldDate = "10/8/2023"
$newDate = "10/9/2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            # msoPlaceholder
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        }
        if (-not $isDatePlaceholder -and $sh.Name -like "Date Placeholder*") {
            $isDatePlaceholder = $true
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout (CustomLayout) under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}
